$d = $word.ActiveDocument

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$pName = $d.Paragraphs.Item(64)
$pWebPage = $d.Paragraphs.Item(68)

$rng = $d.Range($pName.Range.Start, $pWebPage.Range.Start)
$rng.MoveEnd(4, 1)  # wdParagraph = 4
Write-Host ("rng text: [" + $rng.Text + "]")
Write-Host ("rng end: " + $rng.End)
$d.Bookmarks.Add("_GoBack", $rng)
